# Scene.xlsx server-data update
# - Remove the "CloneScene"/Scene2 row and the "RebellerNoob"/SelectScene row
# - Change the RelivePos for the remaining PioneerNoob/villageScene row
# - Change the ID of the remaining (Demo1) row from 4 to 2
# - Move the active selection to F5 (matches the saved worksheet view)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were dropped from the sheet.
# Delete from the bottom up so the earlier row index isn't invalidated.
$ws.Rows(4).Delete()   # was: ../../NFDataCfg/Ini/NFZoneServer/Scene/RebellerNoob/, ID 2, SelectScene
$ws.Rows(2).Delete()   # was: ../../NFDataCfg/Ini/NFZoneServer/Scene/CloneScene/, ID 3, Scene2

# Remaining row 2 (PioneerNoob / villageScene): RelivePos changes from 186,0,88 to 20,0,60
$ws.Cells.Item(2, 5).Value = "20,0,60"

# Remaining row 3 (Demo1 scene): ID changes from 4 to 2
$ws.Cells.Item(3, 2).Value = "2"

# Match the saved selection in the worksheet view
$null = $ws.Range("F5").Select()
